# rebuild feed: fix two typos in existing bios and append three new team
# members (Опанащук, Чернышова, Серикова) to the parsed-data feed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix typo: "дерматовенералогии" -> "дерматовенерологии" (row 2, col C) ---
$cell = $ws.Cells.Item(2, 3)
$text = $cell.Value()
$cell.Value = $text.Replace("дерматовенералогии", "дерматовенерологии")

# --- 2. Fix typo: "Мезфарм" -> "Мезофарм" (row 6, col C) ---
$cell = $ws.Cells.Item(6, 3)
$text = $cell.Value()
$cell.Value = $text.Replace("Мезфарм", "Мезофарм")

# --- 3. Replace placeholder row 13 ("sdfsdf" x3) with real data for
#        Опанащук Марина Николаевна, and add the D/E columns ---
$ws.Cells.Item(13, 1).Value = "https://erabeauty.ru/team/opanashchuk_marina_nikolaevna/"
$ws.Cells.Item(13, 2).Value = "Опанащук Марина Николаевна"
$ws.Cells.Item(13, 3).Value = "Врач-дерматовенеролог, врач-косметолог"
$ws.Cells.Item(13, 4).Value = "2 года"
$ws.Cells.Item(13, 5).Value = "https://erabeauty.ru/wp-content/uploads/2025/04/photo_2025-04-17_14-40-57-768x1024.jpg"

# --- 4. New row 14: Чернышова Мария Михайловна ---
$ws.Cells.Item(14, 1).Value = "https://erabeauty.ru/team/chernyshova_mariya_mihajlovna/"
$ws.Cells.Item(14, 2).Value = "Чернышова Мария Михайловна"
$ws.Cells.Item(14, 3).Value = "Врач-дерматовенеролог, врач-косметолог"
$ws.Cells.Item(14, 4).Value = "2 года"
$ws.Cells.Item(14, 5).Value = "https://erabeauty.ru/wp-content/uploads/2025/04/photo_2025-04-17_14-41-00-768x1024.jpg"

# --- 5. New row 15: Серикова Татьяна Леонидовна ---
$ws.Cells.Item(15, 1).Value = "https://erabeauty.ru/team/serikova_tatyana_leonidovna/"
$ws.Cells.Item(15, 2).Value = "Серикова Татьяна Леонидовна"
$ws.Cells.Item(15, 3).Value = "Медицинская сестра по косметологии, косметолог-эстетист, косметолог-подолог"
$ws.Cells.Item(15, 4).Value = "15 лет"
$ws.Cells.Item(15, 5).Value = "https://erabeauty.ru/wp-content/uploads/2025/04/серикова_татьяна-1-806x1024.jpg"

Write-Output "edit complete"
